# Updates cryptos list prices/volumes (and a few reordered coin rows)
# Mirrors the GitHub Actions scraper re-run: most rows keep their coin/link
# but get a refreshed Price (D) and Volume(1h) (E); a handful of adjacent
# rows swapped rank order so their Coin (B) / Link (C) / Price (D) /
# Volume (E) are all replaced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=rank(unused), B=Coin, C=Link, D=Price, E=Volume(1h)

function Set-Price($row, $price) {
    # Price column holds plain text (e.g. "1.00", "36.305.40") - force text
    # formatting first so Excel doesn't silently coerce it to a Number and
    # drop formatting like trailing zeros.
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $price
}

function Set-Volume($row, $vol) {
    $ws.Cells.Item($row, 5).Value = $vol
}

# Rows where only D (Price) and E (Volume) change
$rows = @(
    @{R=2;  D="36.305.40"; E="  -2.29%  "},
    @{R=3;  D="2.038.77";  E="  -0.85%  "},
    @{R=4;  D="0.999";     E="  -0.06%  "},
    @{R=5;  D="240.20";    E="  -3.39%  "},
    @{R=6;  D="0.663";     E="  -0.46%  "},
    @{R=8;  D="53.49";     E="  -7.99%  "},
    @{R=9;  D="58.13";     E="  -2.95%  "},
    @{R=10; D="0.353";     E="  -8.04%  "},
    @{R=11; D="0.0739";    E="  -5.78%  "},
    @{R=12; D="0.106";     E="  -2.37%  "},
    @{R=15; D="2.335.46";  E="  -0.85%  "},
    @{R=16; D="5.30";      E="  -7.75%  "},
    @{R=17; D="2.073.93";  E="  +0.99%  "},
    @{R=18; D="36.207.73"; E="  -2.49%  "},
    @{R=19; D="16.16";     E="  -14.20%  "},
    @{R=20; D="71.35";     E="  -5.01%  "},
    @{R=21; D="0.0₃0842";  E="  -6.35%  "},
    @{R=24; D="0.999";     E="  -0.12%  "},
    @{R=25; D="2.34";      E="  -6.30%  "},
    @{R=28; D="161.61";    E="  -6.05%  "},
    @{R=29; D="19.99";     E="  -1.67%  "},
    @{R=30; D="0.121";     E="  -3.55%  "},
    @{R=31; D="4.97";      E="  -4.79%  "},
    @{R=32; D="1.13";      E="  -2.71%  "},
    @{R=33; D="4.44";      E="  -5.61%  "},
    @{R=34; D="0.0584";    E="  -7.06%  "},
    @{R=36; D="1.84";      E="  -0.74%  "},
    @{R=37; D="2.17";      E="  -6.53%  "},
    @{R=38; D="0.0806";    E="  -9.41%  "},
    @{R=39; D="1.22";      E="  -9.13%  "},
    @{R=40; D="4.71";      E="  -8.44%  "},
    @{R=41; D="2.83";      E="  -9.05%  "},
    @{R=42; D="0.0211";    E="  -6.59%  "},
    @{R=43; D="1.09";      E="  -6.54%  "},
    @{R=46; D="1.382.24";  E="  +5.51%  "},
    @{R=47; D="15.34";     E="  -11.56%  "},
    @{R=48; D="7.20";      E="  +4.12%  "},
    @{R=51; D="2.221.11";  E="  -1.18%  "}
)

foreach ($item in $rows) {
    Set-Price $item.R $item.D
    Set-Volume $item.R $item.E
}

# Rows where only E (Volume) changes
Set-Volume 7 "  +0.00%  "
Set-Volume 35 "  -0.02%  "

# Rows where the coin (B), link (C), price (D) and volume (E) are all
# replaced (rankings swapped between two adjacent coins)
function Set-CoinRow($row, $coin, $link, $price, $vol) {
    $ws.Cells.Item($row, 2).Value = $coin
    $ws.Cells.Item($row, 3).Value = $link
    Set-Price $row $price
    Set-Volume $row $vol
}

Set-CoinRow 13 "Chainlink"   "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link" "14.53"  "  -10.29%  "
Set-CoinRow 14 "Polygon"     "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"    "0.869"  "  -5.15%  "

Set-CoinRow 22 "BitcoinCash" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch" "235.77" "  -1.40%  "
Set-CoinRow 23 "Uniswap"     "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"         "5.18"   "  -5.17%  "

Set-CoinRow 26 "Cosmos"      "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"     "9.17"   "  -4.61%  "
Set-CoinRow 27 "PancakeSwap" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"    "2.09"   "  -4.28%  "

Set-CoinRow 44 "Aave"        "https://coinranking.com/coin/ixgUfzmLR+aave-aave"           "92.36"  "  -11.08%  "
Set-CoinRow 45 "Cronos"      "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"       "0.0888" "  -11.19%  "

Set-CoinRow 49 "MXToken"     "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"       "2.81"   "  -2.22%  "
Set-CoinRow 50 "RenderToken" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr" "2.26"  "  -7.49%  "
